# Generate Report for Handoff
# Adds a new "dependency" handoff entry (ffff775d73f4-cd37-4502-b88c-38ef1282e163.md)
# to the Overview / zh-cn / de-de sheets, refreshes the hash + handoff timestamps for
# the existing source file (52e84d5b... -> 08c851ac...), and pushes the
# ".localization-config" row down to make room.

$wb = $excel.ActiveWorkbook

$newSrcUuid = "08c851ac-6357-4bb1-a193-b17d1a4cabfd"
$newSrcHash = "2686f8383085294338f2adc279bb23c6fc9f2603"
$newDepFile = "ffff775d73f4-cd37-4502-b88c-38ef1282e163.md"

$zhTime = "2016-02-23 08:00:41"
$deTime = "2016-02-23 08:00:55"

$srcMdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/46aa457e313467aca39eab0e3f0d8830df0670de/e2e/$newSrcUuid.md"
$depMdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/46aa457e313467aca39eab0e3f0d8830df0670de/e2e/$newDepFile"
$configUrl  = "https://github.com/OpenLocalizationTest/oltest/blob/46aa457e313467aca39eab0e3f0d8830df0670de/.localization-config"
$zhXlfName  = "$newSrcUuid.$newSrcHash.zh-cn.xlf"
$deXlfName  = "$newSrcUuid.$newSrcHash.de-de.xlf"
$zhXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cdad7c5e035785c73624033d73fc891e06ab1b32/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/terryjin/ht/$zhXlfName"
$deXlfUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2d3af14b92f964de6e98b11ff9c648ec02b53549/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/terryjin/ht/$deXlfName"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "$newSrcUuid.md"
$ws1.Range("B2").Value = "Ready for handoff"
$ws1.Range("C2").Value = "Ready for handoff"

$ws1.Range("A3").Value = $newDepFile
$ws1.Range("B3").Value = "Ready for handoff"
$ws1.Range("C3").Value = "Ready for handoff"

$ws1.Range("A4").Value = ".localization-config"
$ws1.Range("B4").Value = "Not to be localized"
$ws1.Range("C4").Value = "Not to be localized"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $srcMdUrl, [Type]::Missing, [Type]::Missing, "$newSrcUuid.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), $depMdUrl, [Type]::Missing, [Type]::Missing, $newDepFile)
$ws1.Hyperlinks.Add($ws1.Range("A4"), $configUrl, [Type]::Missing, [Type]::Missing, ".localization-config")

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "$newSrcUuid.md"
$ws2.Range("B2").Value = "Ready for handoff"
$ws2.Range("C2").Value = $zhXlfName
$ws2.Range("D2").Value = $zhTime
$ws2.Range("G2").Value = "0001-01-01 00:00:00"
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = $newDepFile
$ws2.Range("B3").Value = "Ready for handoff"
$ws2.Range("C3").Value = $zhXlfName
$ws2.Range("D3").Value = $zhTime
$ws2.Range("G3").Value = "0001-01-01 00:00:00"
$ws2.Range("H3").Value = "Include"

$ws2.Range("A4").Value = ".localization-config"
$ws2.Range("B4").Value = "Not to be localized"
$ws2.Range("D4").Value = "0001-01-01 00:00:00"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"
$ws2.Range("H4").Value = "Ignored"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $srcMdUrl, [Type]::Missing, [Type]::Missing, "$newSrcUuid.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), $zhXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlfName)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $depMdUrl, [Type]::Missing, [Type]::Missing, $newDepFile)
$ws2.Hyperlinks.Add($ws2.Range("C3"), $zhXlfUrl, [Type]::Missing, [Type]::Missing, $zhXlfName)
$ws2.Hyperlinks.Add($ws2.Range("A4"), $configUrl, [Type]::Missing, [Type]::Missing, ".localization-config")

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = "$newSrcUuid.md"
$ws3.Range("B2").Value = "Ready for handoff"
$ws3.Range("C2").Value = $deXlfName
$ws3.Range("D2").Value = $deTime
$ws3.Range("G2").Value = "0001-01-01 00:00:00"
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = $newDepFile
$ws3.Range("B3").Value = "Ready for handoff"
$ws3.Range("C3").Value = $deXlfName
$ws3.Range("D3").Value = $deTime
$ws3.Range("G3").Value = "0001-01-01 00:00:00"
$ws3.Range("H3").Value = "Include"

$ws3.Range("A4").Value = ".localization-config"
$ws3.Range("B4").Value = "Not to be localized"
$ws3.Range("D4").Value = "0001-01-01 00:00:00"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"
$ws3.Range("H4").Value = "Ignored"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $srcMdUrl, [Type]::Missing, [Type]::Missing, "$newSrcUuid.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), $deXlfUrl, [Type]::Missing, [Type]::Missing, $deXlfName)
$ws3.Hyperlinks.Add($ws3.Range("A3"), $depMdUrl, [Type]::Missing, [Type]::Missing, $newDepFile)
$ws3.Hyperlinks.Add($ws3.Range("C3"), $deXlfUrl, [Type]::Missing, [Type]::Missing, $deXlfName)
$ws3.Hyperlinks.Add($ws3.Range("A4"), $configUrl, [Type]::Missing, [Type]::Missing, ".localization-config")

$wb.Save()
